$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) followed by 10 data rows (rows 2-11).
# Row 2 (siteid 1771088 / recording_id 26834047) is removed entirely.
# Of the remaining 9 rows, the block that is now at rows 7-11 (5 rows, the
# "2020/2021" records) needs to move above the block now at rows 2-6 that
# remains (4 rows, the older "2018" records), preserving each block's
# internal order.

# Step 1: delete row 2 - remaining rows shift up so the sheet now has 9 data rows (2-10).
#   rows 2-5  -> older block (4 rows)
#   rows 6-10 -> newer block (5 rows)
$ws.Rows.Item(2).Delete()

# Step 2: make room for the newer block above the older block by inserting
# 5 blank rows at rows 2-6 (only within the used columns A:K so we do not
# touch unrelated columns on the row).
$ws.Range("A2:K6").Insert()

# After the insert:
#   rows 2-6   -> blank
#   rows 7-10  -> older block (4 rows)
#   rows 11-15 -> newer block (5 rows)

# Step 3: move (cut/paste) the newer block into the blank rows at the top.
$ws.Range("A11:K15").Cut($ws.Range("A2"))

# Step 4: remove the now-empty leftover rows at the bottom.
$ws.Range("A11:K15").Delete()
